$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10, shifting existing rows 10:79 down to 11:80,
# carrying formatting from the row above (default Excel behaviour).
$ws.Rows("10:10").Insert(1)

# Populate the newly inserted row 10 with the new record values.
$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(10, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(10, 3).Value = "La Araucanía"
$ws.Cells.Item(10, 4).Value = 44881
$ws.Cells.Item(10, 5).Value = 9
$ws.Cells.Item(10, 6).Value = 300000000
$ws.Cells.Item(10, 7).Value = "Espárragos"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 200
$ws.Cells.Item(10, 11).Value = 1500
$ws.Cells.Item(10, 12).Value = 1500
$ws.Cells.Item(10, 13).Value = 1500
$ws.Cells.Item(10, 14).Value = "`$/kilo"
$ws.Cells.Item(10, 15).Value = "Región del Maule"
$ws.Cells.Item(10, 16).Value = 1500
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = "Hortaliza"
